$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(42613.76190972222, 56),
    @(42613.891898148147, 17),
    @(42614.889270833337, 19),
    @(42615.887835648151, 96)
)

$row = 8
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    for ($col = 3; $col -le 13; $col++) {
        $ws.Cells.Item($row, $col).Value = 0
    }
    $ws.Cells.Item($row, 14).Value = "Random"
    $row++
}
